$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped on Sat Sep 30 09:17:46 UTC 2023
$ws.Range("D2").Value = "26.968.24"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "1.677.21"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.00"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.34"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D12").Value = "1.912.67"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "1.685.04"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.86"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "26.979.25"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.35"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +4.00%  "
$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.21"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.51"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  +1.48%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.33"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "1.487.15"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0174"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.900"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  -3.36%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +2.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.51"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "1.816.70"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.778"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.53"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +15.56%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.53"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0510"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.46%  "
